$wb = $excel.ActiveWorkbook

$oldGuid = "a0066752-6342-4d87-87f5-4ccaa044e04b"
$newGuid = "5d75d24b-742c-475e-a68c-7373f9cb30e9"
$oldHash = "61020bdb38a25ee37e2ff3c68d191e82965e658a"
$newHash = "abfbfbe38288b283a26682d9f275516bb587df2f"

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value = "$newGuid.md"
$wsOverview.Range("D2").Value = "2016-49-19 16:49:00"

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("A2").Value = "$newGuid.md"
$wsZhCn.Range("D2").Value = "$newGuid.$newHash.zh-cn.xlf"
$wsZhCn.Range("E2").Value = "2016-03-19 16:48:57"

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("A2").Value = "$newGuid.md"
$wsDeDe.Range("D2").Value = "$newGuid.$newHash.de-de.xlf"
$wsDeDe.Range("E2").Value = "2016-03-19 16:49:00"
